$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 43,5
$data[0,0] = "Rv3418c"
$data[0,1] = 5
$data[0,2] = "groS cpn10 groES mopB Rv3418c MTCY78.11"
$data[0,3] = "FUNCTION: Binds to Cpn60 in the presence of Mg-ATP and suppresses the ATPase activity of the latter."
$data[0,4] = 36
$data[1,0] = "Rv3247c"
$data[1,1] = 5
$data[1,2] = "tmk Rv3247c"
$data[1,3] = "FUNCTION: Catalyzes the reversible phosphorylation of deoxythymidine monophosphate (dTMP) to deoxythymidine diphosphate (dTDP), using ATP as its preferred phosphoryl donor. Situated at the junction of both de novo and salvage pathways of deoxythymidine triphosphate (dTTP) synthesis, is essential for DNA synthesis and cellular growth. Has a broad specificity for nucleoside triphosphates, being highly active with ATP or dATP as phosphate donors, and less active with ITP, GTP, CTP and UTP."
$data[1,4] = 36
$data[2,0] = "Rv3280"
$data[2,1] = 5
$data[2,2] = "accD5 pccB Rv3280 MTCY71.20"
$data[2,3] = "FUNCTION: Component of a biotin-dependent acyl-CoA carboxylase complex. This subunit transfers the CO2 from carboxybiotin to the CoA ester substrate (PubMed:16354663, PubMed:16385038, PubMed:28222482). When associated with the alpha3 subunit AccA3, is involved in the carboxylation of acetyl-CoA and propionyl-CoA, with a preference for propionyl-CoA (PubMed:16354663, PubMed:16385038, PubMed:28222482). Is also required for the activity of the long-chain acyl-CoA carboxylase (LCC) complex (PubMed:28222482). {ECO:0000269|PubMed:16354663, ECO:0000269|PubMed:16385038, ECO:0000269|PubMed:28222482}."
$data[2,4] = 36
$data[3,0] = "Rv1617"
$data[3,1] = 5
$data[3,2] = "pyk pykA Rv1617 MTCY01B2.09"
$data[3,3] = ""
$data[3,4] = 36
$data[4,0] = "Rv1093"
$data[4,1] = 5
$data[4,2] = "glyA1 glyA Rv1093 MTV017.46"
$data[4,3] = "FUNCTION: Catalyzes the reversible interconversion of serine and glycine with tetrahydrofolate (THF) serving as the one-carbon carrier. This reaction serves as the major source of one-carbon groups required for the biosynthesis of purines, thymidylate, methionine, and other important biomolecules. Also exhibits THF-independent aldolase activity toward beta-hydroxyamino acids, producing glycine and aldehydes, via a retro-aldol mechanism. Thus, is able to catalyze the cleavage of L-allo-threonine. {ECO:0000269|PubMed:12913008}."
$data[4,4] = 36
$data[5,0] = "Rv3246c"
$data[5,1] = 5
$data[5,2] = "mtrA Rv3246c MTCY20B11.21c"
$data[5,3] = "FUNCTION: Member of the two-component regulatory system MtrA/MtrB. Binds direct repeat motifs of sequence 5'-GTCACAGCG-3', phosphorylation confers higher affinity. Overexpression decreases bacteria viability upon infection of human THP-1 macrophage cell line, due at least in part to impaired blockage of phagosome-lysosome fusion (upon infection bacteria usually remain in phagosomes). Infecting C57BL/6 mice with an overexpressing strain leads to an attentuated infection in both spleen and lungs. The level of dnaA mRNA increases dramatically. Binds the promoter of dnaA, fbpD, ripA and itself, as well as oriC, which it may regulate. Upon co-overexpression of MrtA and MtrB growth in macrophages is partially restored, dnaA expression is not induced, although mouse infections are still attenuated, suggesting that bacterial growth in macrophages requires an optimal ratio of MtrB to MtrA. {ECO:0000269|PubMed:20223818, ECO:0000269|PubMed:21295603, ECO:0000269|PubMed:22610443}."
$data[5,4] = 36
$data[6,0] = "Rv2747"
$data[6,1] = 5
$data[6,2] = "argA Rv2747"
$data[6,3] = "FUNCTION: Catalyzes the conversion of L-glutamate to alpha-N-acetyl-L-glutamate. L-glutamine is a significantly better substrate compared to L-glutamate. {ECO:0000269|PubMed:15838030}."
$data[6,4] = 36
$data[7,0] = "Rv3219"
$data[7,1] = 5
$data[7,2] = "whiB1 Rv3219"
$data[7,3] = "FUNCTION: Acts as a transcriptional repressor, inhibiting expression in vitro. Probably redox-responsive. The apo- but not holo-form binds to its own promoter as well as that of groEL2. Oxidized apo-form and nitrosylated holo-form also bind DNA. The apo-form has been shown to act as a protein disulfide reductase (PubMed:17157031) (PubMed:19016840), but also not to act as a protein disulfide reductase (PubMed:20929442). {ECO:0000269|PubMed:17157031, ECO:0000269|PubMed:19016840, ECO:0000269|PubMed:20929442, ECO:0000269|PubMed:22464736}."
$data[7,4] = 36
$data[8,0] = "Rv3372"
$data[8,1] = 5
$data[8,2] = "otsB otsB2 Rv3372"
$data[8,3] = "FUNCTION: Removes the phosphate from trehalose 6-phosphate to produce free trehalose. {ECO:0000269|PubMed:15158675, ECO:0000269|PubMed:15703182}."
$data[8,4] = 36
$data[9,0] = "Rv3285"
$data[9,1] = 5
$data[9,2] = "accA3 Rv3285"
$data[9,3] = "FUNCTION: Component of a biotin-dependent acyl-CoA carboxylase complex. This subunit catalyzes the ATP-dependent carboxylation of the biotin carried by the biotin carboxyl carrier (BCC) domain, resulting in the formation of carboxyl biotin (PubMed:16354663, PubMed:16385038, PubMed:17114269). When associated with the beta5 subunit AccD5, is involved in the carboxylation of acetyl-CoA and propionyl-CoA, with a preference for propionyl-CoA (PubMed:16354663, PubMed:16385038). When associated with the beta6 subunit AccD6, is involved in the carboxylation of acetyl-CoA and propionyl-CoA, with a preference for acetyl-CoA (PubMed:17114269). When associated with the beta4 subunit AccD4, the beta5 subunit AccD5 and the epsilon subunit AccE5, forms the LCC complex, which is involved in the carboxylation of long chain acyl-CoA (PubMed:16354663, PubMed:28222482). The LCC complex can use C16-C24 substrates, the highest specific activity is obtained with carboxy-C20-CoA (PubMed:28222482). {ECO:0000269|PubMed:16354663, ECO:0000269|PubMed:16385038, ECO:0000269|PubMed:17114269, ECO:0000269|PubMed:28222482}."
$data[9,4] = 36
$data[10,0] = "Rv3042c"
$data[10,1] = 5
$data[10,2] = "serB2 Rv3042c"
$data[10,3] = "FUNCTION: Catalyzes the dephosphorylation of O-phospho-L-serine into L-serine, a step in the L-serine biosynthetic pathway (PubMed:25037224, PubMed:25521849). Exhibits high specificity for L-phosphoserine compared to substrates like L-phosphothreonine (5% relative activity) and L-phosphotyrosine (1.7% relative activity) (PubMed:25521849). {ECO:0000269|PubMed:25037224, ECO:0000269|PubMed:25521849}.; FUNCTION: In the host, induces significant cytoskeleton rearrangements through cofilin dephosphorylation and its subsequent activation, and affects the expression of genes that regulate actin dynamics. It specifically interacts with HSP90, HSP70 and HSP27 that block apoptotic pathways but not with other HSPs. Also interacts with GAPDH. It actively dephosphorylates MAP kinase p38 and NF-kappa B p65 (specifically at Ser-536) that play crucial roles in inflammatory and immune responses. This in turn leads to down-regulation of Interleukin 8, a chemotactic and inflammatory cytokine. Thus might help the pathogen to evade the host's immune response (PubMed:26984196). Exogenous addition of purified SerB2 protein to human THP-1 cells (that can be differentiated into macrophage-like cells) induces microtubule rearrangements; the phosphatase activity is co-related to the elicited rearrangements, while addition of the ACT-domains alone elicits no rearrangements (PubMed:25521849). {ECO:0000269|PubMed:25521849, ECO:0000269|PubMed:26984196}."
$data[10,4] = 36
$data[11,0] = "Rv2754c"
$data[11,1] = 5
$data[11,2] = "thyX Rv2754c MTV002.19c"
$data[11,3] = "FUNCTION: Catalyzes the reductive methylation of 2'-deoxyuridine-5'-monophosphate (dUMP) to 2'-deoxythymidine-5'-monophosphate (dTMP) while utilizing 5,10-methylenetetrahydrofolate (mTHF) as the methyl donor, and NADPH and FADH(2) as the reductant (PubMed:18493582). Is essential for growth of the pathogen on solid media in vitro; the essential function is something other than dTMP synthase (PubMed:12657046) (PubMed:22034487). {ECO:0000269|PubMed:12657046, ECO:0000269|PubMed:16139296, ECO:0000269|PubMed:18493582, ECO:0000269|PubMed:22034487}."
$data[11,4] = 36
$data[12,0] = "Rv3260c"
$data[12,1] = 5
$data[12,2] = "whiB2 Rv3260c"
$data[12,3] = "FUNCTION: Acts as a transcriptional regulator. Probably redox-responsive. The apo- but not holo-form probably binds DNA (By similarity). {ECO:0000250}.; FUNCTION: The apo-form functions as a chaperone, preventing aggregation or helping in correct refolding of a number of substrates; this activity does not require ATP or the ability to bind a Fe-S cluster. Chaperone activity is insensitive to the redox state of its cysteine residues. The apo-form has no protein disulfide reductase activity. The apo-form binds to its own promoter. {ECO:0000269|PubMed:19016840, ECO:0000269|PubMed:22686939}."
$data[12,4] = 36
$data[13,0] = "Rv3245c"
$data[13,1] = 5
$data[13,2] = "mtrB Rv3245c MTCY20B11.20c"
$data[13,3] = "FUNCTION: Member of the two-component regulatory system MtrA/MtrB. Probably functions as a membrane-associated protein kinase that phosphorylates MtrA in response to environmental signals. Autophosphorylates and transfers phosphate to MtrA in vitro. Overexpression of MtrA alone decreases bacterial virulence in mouse infection; co-expression of MtrA and MtrB restores normal bacterial growth, suggesting that bacterial growth in macrophages requires an optimal ratio of MtrB to MtrA. Probably plays a role in cell division. {ECO:0000269|PubMed:21295603, ECO:0000269|PubMed:22610443}."
$data[13,4] = 36
$data[14,0] = "Rv3789"
$data[14,1] = 4
$data[14,2] = "Rv3789 MTCY13D12.23"
$data[14,3] = "FUNCTION: Required for arabinosylation of arabinogalactan (AG), an essential component of the mycobacterial cell wall. Probably acts as an anchor protein recruiting AftA, the first arabinosyl transferase involved in AG biosynthesis. {ECO:0000269|PubMed:26369580}."
$data[14,4] = 36
$data[15,0] = "Rv2477c"
$data[15,1] = 4
$data[15,2] = "ettA Rv2477c"
$data[15,3] = "FUNCTION: A translation factor that gates the progression of the 70S ribosomal initiation complex (IC, containing tRNA(fMet) in the P-site) into the translation elongation cycle by using a mechanism sensitive to the ATP/ADP ratio. Binds to the 70S ribosome E-site where it modulates the state of the translating ribosome during subunit translocation. ATP hydrolysis probably frees it from the ribosome, which can enter the elongation phase. {ECO:0000255|HAMAP-Rule:MF_00847}."
$data[15,4] = 36
$data[16,0] = "Rv3266c"
$data[16,1] = 4
$data[16,2] = "rmlD Rv3266c"
$data[16,3] = "FUNCTION: Involved in the biosynthesis of the dTDP-L-rhamnose which is a component of the critical linker, D-N-acetylglucosamine-L-rhamnose disaccharide, which connects the galactan region of arabinogalactan to peptidoglycan via a phosphodiester linkage (PubMed:12029057). Catalyzes the reduction of dTDP-6-deoxy-L-lyxo-4-hexulose to yield dTDP-L-rhamnose (By similarity). {ECO:0000250|UniProtKB:P26392, ECO:0000269|PubMed:12029057}."
$data[16,4] = 36
$data[17,0] = "Rv3215"
$data[17,1] = 3
$data[17,2] = "menF entC Rv3215"
$data[17,3] = "FUNCTION: Catalyzes the conversion of chorismate to isochorismate. {ECO:0000250|UniProtKB:P38051}."
$data[17,4] = 36
$data[18,0] = "Rv0337c"
$data[18,1] = 3
$data[18,2] = "aspC Rv0337c MTCY279.04c"
$data[18,3] = ""
$data[18,4] = 36
$data[19,0] = "Rv1224"
$data[19,1] = 3
$data[19,2] = "tatB Rv1224 MTCI61.07"
$data[19,3] = "FUNCTION: Part of the twin-arginine translocation (Tat) system that transports large folded proteins containing a characteristic twin-arginine motif in their signal peptide across membranes. Together with TatC, TatB is part of a receptor directly interacting with Tat signal peptides. TatB may form an oligomeric binding site that transiently accommodates folded Tat precursor proteins before their translocation. {ECO:0000255|HAMAP-Rule:MF_00237}."
$data[19,4] = 36
$data[20,0] = "Rv3255c"
$data[20,1] = 3
$data[20,2] = "manA Rv3255c"
$data[20,3] = ""
$data[20,4] = 36
$data[21,0] = "Rv1338"
$data[21,1] = 3
$data[21,2] = "murI Rv1338 MTCY02B10.02 MTCY130.23"
$data[21,3] = "FUNCTION: Provides the (R)-glutamate required for cell wall biosynthesis. {ECO:0000255|HAMAP-Rule:MF_00258}."
$data[21,4] = 36
$data[22,0] = "Rv0041"
$data[22,1] = 3
$data[22,2] = "leuS Rv0041 MTCY21D4.04"
$data[22,3] = ""
$data[22,4] = 36
$data[23,0] = "Rv3627c"
$data[23,1] = 3
$data[23,2] = "Rv3627c"
$data[23,3] = "FUNCTION: Carboxypeptidase that cleaves terminal D-alanine from peptidoglycan in the mycobacterial cell wall. May cleave L-Lys-D-Ala and/or D-Ala-D-Ala peptide bonds. Exerts important effects on mycobacterial cell morphology and cell division. {ECO:0000269|PubMed:31000162}."
$data[23,4] = 36
$data[24,0] = "Rv2555c"
$data[24,1] = 3
$data[24,2] = "alaS Rv2555c MTCY159.01 MTCYW318.01c"
$data[24,3] = "FUNCTION: Catalyzes the attachment of alanine to tRNA(Ala) in a two-step reaction: alanine is first activated by ATP to form Ala-AMP and then transferred to the acceptor end of tRNA(Ala). Also edits incorrectly charged Ser-tRNA(Ala) and Gly-tRNA(Ala) via its editing domain. {ECO:0000255|HAMAP-Rule:MF_00036}."
$data[24,4] = 36
$data[25,0] = "Rv3598c"
$data[25,1] = 3
$data[25,2] = "lysS1 lysS Rv3598c MTCY07H7B.24"
$data[25,3] = ""
$data[25,4] = 36
$data[26,0] = "Rv1402"
$data[26,1] = 3
$data[26,2] = "priA Rv1402 MTCY21B4.19"
$data[26,3] = "FUNCTION: Involved in the restart of stalled replication forks. Recognizes and binds the arrested nascent DNA chain at stalled replication forks. It can open the DNA duplex, via its helicase activity, and promote assembly of the primosome and loading of the major replicative helicase DnaB onto DNA. {ECO:0000255|HAMAP-Rule:MF_00983}."
$data[26,4] = 36
$data[27,0] = "Rv0482"
$data[27,1] = 3
$data[27,2] = "murB Rv0482 MTCY20G9.08"
$data[27,3] = "FUNCTION: Cell wall formation. {ECO:0000250}."
$data[27,4] = 36
$data[28,0] = "Rv0638"
$data[28,1] = 3
$data[28,2] = "secE Rv0638 MTCY20H10.19"
$data[28,3] = "FUNCTION: Essential subunit of the Sec protein translocation channel SecYEG. Clamps together the 2 halves of SecY. May contact the channel plug during translocation. {ECO:0000255|HAMAP-Rule:MF_00422}."
$data[28,4] = 36
$data[29,0] = "Rv0285"
$data[29,1] = 3
$data[29,2] = "PE5 Rv0285 LH57_01560"
$data[29,3] = "FUNCTION: Important for the siderophore-mediated iron-acquisition function of ESX-3 (PubMed:26729876). May play a pivotal role in the evasion of host immune response by M.tuberculosis. Mediates production of IL-10 via activation of the p38 and ERK1/2 mitogen-activated protein kinase (MAPK) signaling pathways (PubMed:23284742). {ECO:0000269|PubMed:23284742, ECO:0000269|PubMed:26729876}."
$data[29,4] = 36
$data[30,0] = "Rv1828"
$data[30,1] = 3
$data[30,2] = "Rv1828 MTCY1A11.15c"
$data[30,3] = "FUNCTION: Transcriptional regulator that binds to its own promoter and thus may play a role in the regulation of the cotranscribed genes Rv1827 and Rv1828. Can also bind several promoter regions of genes that are essential, including ftsZ. Binds to the imperfect everted repeat sequence CTCAA through its winged-HTH motif. {ECO:0000269|PubMed:30306715}."
$data[30,4] = 36
$data[31,0] = "Rv2093c"
$data[31,1] = 3
$data[31,2] = "tatC Rv2093c MTCY49.33c"
$data[31,3] = "FUNCTION: Part of the twin-arginine translocation (Tat) system that transports large folded proteins containing a characteristic twin-arginine motif in their signal peptide across membranes. Together with TatB, TatC is part of a receptor directly interacting with Tat signal peptides. {ECO:0000255|HAMAP-Rule:MF_00902}."
$data[31,4] = 36
$data[32,0] = "Rv3336c"
$data[32,1] = 3
$data[32,2] = "trpS Rv3336c MTV016.36c"
$data[32,3] = "FUNCTION: Catalyzes the attachment of tryptophan to tRNA(Trp). {ECO:0000255|HAMAP-Rule:MF_00140}."
$data[32,4] = 36
$data[33,0] = "Rv0384c"
$data[33,1] = 3
$data[33,2] = "clpB Rv0384c MTV036.19c"
$data[33,3] = "FUNCTION: Part of a stress-induced multi-chaperone system, it is involved in the recovery of the cell from heat-induced damage, in cooperation with DnaK, DnaJ and GrpE. Acts before DnaK, in the processing of protein aggregates. Protein binding stimulates the ATPase activity; ATP hydrolysis unfolds the denatured protein aggregates, which probably helps expose new hydrophobic binding sites on the surface of ClpB-bound aggregates, contributing to the solubilization and refolding of denatured protein aggregates by DnaK (By similarity). {ECO:0000250}."
$data[33,4] = 36
$data[34,0] = "Rv1296"
$data[34,1] = 2
$data[34,2] = "thrB Rv1296 MTCY373.16"
$data[34,3] = "FUNCTION: Catalyzes the ATP-dependent phosphorylation of L-homoserine to L-homoserine phosphate. {ECO:0000255|HAMAP-Rule:MF_00384}."
$data[34,4] = 36
$data[35,0] = "Rv3302c"
$data[35,1] = 2
$data[35,2] = "glpD2 Rv3302c MTCI418A.04c MTV016.01c"
$data[35,3] = ""
$data[35,4] = 36
$data[36,0] = "Rv1292"
$data[36,1] = 2
$data[36,2] = "argS Rv1292 MTCY373.12"
$data[36,3] = ""
$data[36,4] = 36
$data[37,0] = "Rv3221c"
$data[37,1] = 2
$data[37,2] = "Rv3221c MTCY07D11.05"
$data[37,3] = ""
$data[37,4] = 36
$data[38,0] = "Rv3053c"
$data[38,1] = 2
$data[38,2] = "nrdH Rv3053c"
$data[38,3] = "FUNCTION: Electron transport system for the ribonucleotide reductase system NrdEF. {ECO:0000256|ARBA:ARBA00002292}."
$data[38,4] = 36
$data[39,0] = "Rv2603c"
$data[39,1] = 1
$data[39,2] = "Rv2603c MTCI270A.02"
$data[39,3] = ""
$data[39,4] = 36
$data[40,0] = "Rv3264c"
$data[40,1] = 1
$data[40,2] = "manB Rv3264c"
$data[40,3] = ""
$data[40,4] = 36
$data[41,0] = "Rv3222c"
$data[41,1] = 1
$data[41,2] = "Rv3222c"
$data[41,3] = ""
$data[41,4] = 36
$data[42,0] = "Rv1122"
$data[42,1] = 1
$data[42,2] = "gnd2 Rv1122"
$data[42,3] = ""
$data[42,4] = 36

$ws.Range("A2:E44").Value = $data
